$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = '@'
    $cell.Value = $val
    $cell.Style = 'Normal'
}

function Set-HeaderCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# ---- Sheet: Overall ----
$ws = $wb.Worksheets.Item('Overall')

# Header row
Set-HeaderCell $ws 1 1 'Share of 990 filers with government grants at risk'
Set-HeaderCell $ws 1 2 'Number of 990 filers with government grants'
Set-HeaderCell $ws 1 3 'Total government grants ($)'
Set-HeaderCell $ws 1 4 'Size of operating surplus with government grants'
Set-HeaderCell $ws 1 5 'Size of operating surplus without government grants'

# Data rows
Set-TextCell $ws 2 1 '70.00%'
Set-TextCell $ws 2 2 '630'
Set-TextCell $ws 2 3 '$1,081,113,515'
Set-TextCell $ws 2 4 '8.92%'
Set-TextCell $ws 2 5 '-17.37%'


# ---- Sheet: County ----
$ws = $wb.Worksheets.Item('County')

# Header row
Set-HeaderCell $ws 1 1 'Geography'
Set-HeaderCell $ws 1 2 'Share of 990 filers with government grants at risk'
Set-HeaderCell $ws 1 3 'Number of 990 filers with government grants'
Set-HeaderCell $ws 1 4 'Total government grants ($)'
Set-HeaderCell $ws 1 5 'Size of operating surplus with government grants'
Set-HeaderCell $ws 1 6 'Size of operating surplus without government grants'

# Data rows
Set-TextCell $ws 2 1 'United States'
Set-TextCell $ws 2 2 '67.35%'
Set-TextCell $ws 2 3 '103,475'
Set-TextCell $ws 2 4 '$267,700,640,005'
Set-TextCell $ws 2 5 '9.05%'
Set-TextCell $ws 2 6 '-12.83%'

Set-TextCell $ws 3 1 'Hawaii'
Set-TextCell $ws 3 2 '70.00%'
Set-TextCell $ws 3 3 '630'
Set-TextCell $ws 3 4 '$1,081,113,515'
Set-TextCell $ws 3 5 '8.92%'
Set-TextCell $ws 3 6 '-17.37%'

Set-TextCell $ws 4 1 'Hawaii County'
Set-TextCell $ws 4 2 '73.15%'
Set-TextCell $ws 4 3 '108'
Set-TextCell $ws 4 4 '$93,231,648'
Set-TextCell $ws 4 5 '6.89%'
Set-TextCell $ws 4 6 '-16.25%'

Set-TextCell $ws 5 1 'Honolulu County'
Set-TextCell $ws 5 2 '71.78%'
Set-TextCell $ws 5 3 '365'
Set-TextCell $ws 5 4 '$874,969,518'
Set-TextCell $ws 5 5 '8.42%'
Set-TextCell $ws 5 6 '-19.14%'

Set-TextCell $ws 6 1 'Kauai County'
Set-TextCell $ws 6 2 '52.50%'
Set-TextCell $ws 6 3 '40'
Set-TextCell $ws 6 4 '$13,473,937'
Set-TextCell $ws 6 5 '20.31%'
Set-TextCell $ws 6 6 '-0.86%'

Set-TextCell $ws 7 1 'Maui County'
Set-TextCell $ws 7 2 '67.52%'
Set-TextCell $ws 7 3 '117'
Set-TextCell $ws 7 4 '$99,438,412'
Set-TextCell $ws 7 5 '11.26%'
Set-TextCell $ws 7 6 '-18.60%'


# ---- Sheet: Congressional District ----
$ws = $wb.Worksheets.Item('Congressional District')

# Header row
Set-HeaderCell $ws 1 1 'Geography'
Set-HeaderCell $ws 1 2 'Share of 990 filers with government grants at risk'
Set-HeaderCell $ws 1 3 'Number of 990 filers with government grants'
Set-HeaderCell $ws 1 4 'Total government grants ($)'
Set-HeaderCell $ws 1 5 'Size of operating surplus with government grants'
Set-HeaderCell $ws 1 6 'Size of operating surplus without government grants'

# Data rows
Set-TextCell $ws 2 1 'United States'
Set-TextCell $ws 2 2 '67.35%'
Set-TextCell $ws 2 3 '103,475'
Set-TextCell $ws 2 4 '$267,700,640,005'
Set-TextCell $ws 2 5 '9.05%'
Set-TextCell $ws 2 6 '-12.83%'

Set-TextCell $ws 3 1 'Hawaii'
Set-TextCell $ws 3 2 '70.00%'
Set-TextCell $ws 3 3 '630'
Set-TextCell $ws 3 4 '$1,081,113,515'
Set-TextCell $ws 3 5 '8.92%'
Set-TextCell $ws 3 6 '-17.37%'

Set-TextCell $ws 4 1 'Congressional District 1'
Set-TextCell $ws 4 2 '71.48%'
Set-TextCell $ws 4 3 '298'
Set-TextCell $ws 4 4 '$790,093,863'
Set-TextCell $ws 4 5 '8.27%'
Set-TextCell $ws 4 6 '-19.57%'

Set-TextCell $ws 5 1 'Congressional District 2'
Set-TextCell $ws 5 2 '68.67%'
Set-TextCell $ws 5 3 '332'
Set-TextCell $ws 5 4 '$291,019,652'
Set-TextCell $ws 5 5 '9.73%'
Set-TextCell $ws 5 6 '-15.63%'


# ---- Sheet: Size ----
$ws = $wb.Worksheets.Item('Size')

# Header row
Set-HeaderCell $ws 1 1 'Size'
Set-HeaderCell $ws 1 2 'Share of 990 filers with government grants at risk'
Set-HeaderCell $ws 1 3 'Number of 990 filers with government grants'
Set-HeaderCell $ws 1 4 'Total government grants ($)'
Set-HeaderCell $ws 1 5 'Size of operating surplus with government grants'
Set-HeaderCell $ws 1 6 'Size of operating surplus without government grants'

# Data rows
Set-TextCell $ws 2 1 'Between $100K and $499K'
Set-TextCell $ws 2 2 '70.15%'
Set-TextCell $ws 2 3 '201'
Set-TextCell $ws 2 4 '$24,503,478'
Set-TextCell $ws 2 5 '11.95%'
Set-TextCell $ws 2 6 '-18.22%'

Set-TextCell $ws 3 1 'Between $1M and $4.99M'
Set-TextCell $ws 3 2 '72.00%'
Set-TextCell $ws 3 3 '175'
Set-TextCell $ws 3 4 '$171,732,848'
Set-TextCell $ws 3 5 '9.03%'
Set-TextCell $ws 3 6 '-19.27%'

Set-TextCell $ws 4 1 'Between $500K and $999K'
Set-TextCell $ws 4 2 '70.09%'
Set-TextCell $ws 4 3 '117'
Set-TextCell $ws 4 4 '$31,929,560'
Set-TextCell $ws 4 5 '8.80%'
Set-TextCell $ws 4 6 '-14.41%'

Set-TextCell $ws 5 1 'Between $5M and $9.99M'
Set-TextCell $ws 5 2 '69.05%'
Set-TextCell $ws 5 3 '42'
Set-TextCell $ws 5 4 '$132,115,397'
Set-TextCell $ws 5 5 '5.52%'
Set-TextCell $ws 5 6 '-28.85%'

Set-TextCell $ws 6 1 'Greater than $10M'
Set-TextCell $ws 6 2 '62.69%'
Set-TextCell $ws 6 3 '67'
Set-TextCell $ws 6 4 '$715,119,205'
Set-TextCell $ws 6 5 '5.73%'
Set-TextCell $ws 6 6 '-7.62%'

Set-TextCell $ws 7 1 'Less than $100K'
Set-TextCell $ws 7 2 '75.00%'
Set-TextCell $ws 7 3 '28'
Set-TextCell $ws 7 4 '$5,713,027'
Set-TextCell $ws 7 5 '17.69%'
Set-TextCell $ws 7 6 '-24.27%'

Set-TextCell $ws 8 1 'Total'
Set-TextCell $ws 8 2 '70.00%'
Set-TextCell $ws 8 3 '630'
Set-TextCell $ws 8 4 '$1,081,113,515'
Set-TextCell $ws 8 5 '8.92%'
Set-TextCell $ws 8 6 '-17.37%'


# ---- Sheet: Subsector ----
$ws = $wb.Worksheets.Item('Subsector')

# Header row
Set-HeaderCell $ws 1 1 'Subsector'
Set-HeaderCell $ws 1 2 'Share of 990 filers with government grants at risk'
Set-HeaderCell $ws 1 3 'Number of 990 filers with government grants'
Set-HeaderCell $ws 1 4 'Total government grants ($)'
Set-HeaderCell $ws 1 5 'Size of operating surplus with government grants'
Set-HeaderCell $ws 1 6 'Size of operating surplus without government grants'

# Data rows
Set-TextCell $ws 2 1 'Arts, Culture, and Humanities'
Set-TextCell $ws 2 2 '67.95%'
Set-TextCell $ws 2 3 '78'
Set-TextCell $ws 2 4 '$52,913,318'
Set-TextCell $ws 2 5 '10.85%'
Set-TextCell $ws 2 6 '-16.11%'

Set-TextCell $ws 3 1 'Education (Excluding Universities)'
Set-TextCell $ws 3 2 '70.77%'
Set-TextCell $ws 3 3 '65'
Set-TextCell $ws 3 4 '$34,850,435'
Set-TextCell $ws 3 5 '9.14%'
Set-TextCell $ws 3 6 '-12.65%'

Set-TextCell $ws 4 1 'Environment and Animals'
Set-TextCell $ws 4 2 '65.28%'
Set-TextCell $ws 4 3 '72'
Set-TextCell $ws 4 4 '$39,214,536'
Set-TextCell $ws 4 5 '13.14%'
Set-TextCell $ws 4 6 '-17.78%'

Set-TextCell $ws 5 1 'Health (Excluding Hospitals)'
Set-TextCell $ws 5 2 '64.52%'
Set-TextCell $ws 5 3 '62'
Set-TextCell $ws 5 4 '$144,808,995'
Set-TextCell $ws 5 5 '7.27%'
Set-TextCell $ws 5 6 '-20.35%'

Set-TextCell $ws 6 1 'Hospitals'
Set-TextCell $ws 6 2 '42.86%'
Set-TextCell $ws 6 3 '7'
Set-TextCell $ws 6 4 '$77,923,568'
Set-TextCell $ws 6 5 '9.43%'
Set-TextCell $ws 6 6 '0.32%'

Set-TextCell $ws 7 1 'Human Services'
Set-TextCell $ws 7 2 '69.03%'
Set-TextCell $ws 7 3 '155'
Set-TextCell $ws 7 4 '$247,080,925'
Set-TextCell $ws 7 5 '10.81%'
Set-TextCell $ws 7 6 '-18.43%'

Set-TextCell $ws 8 1 'International, Foreign Affairs'
Set-TextCell $ws 8 2 '80.00%'
Set-TextCell $ws 8 3 '5'
Set-TextCell $ws 8 4 '$833,746'
Set-TextCell $ws 8 5 '20.57%'
Set-TextCell $ws 8 6 '-14.41%'

Set-TextCell $ws 9 1 'Public, Societal Benefit'
Set-TextCell $ws 9 2 '65.45%'
Set-TextCell $ws 9 3 '55'
Set-TextCell $ws 9 4 '$212,779,190'
Set-TextCell $ws 9 5 '5.10%'
Set-TextCell $ws 9 6 '-15.30%'

Set-TextCell $ws 10 1 'Religion Related'
Set-TextCell $ws 10 2 '66.67%'
Set-TextCell $ws 10 3 '3'
Set-TextCell $ws 10 4 '$12,634,018'
Set-TextCell $ws 10 5 '-3.97%'
Set-TextCell $ws 10 6 '-30.02%'

Set-TextCell $ws 11 1 'Unclassified'
Set-TextCell $ws 11 2 '80.49%'
Set-TextCell $ws 11 3 '123'
Set-TextCell $ws 11 4 '$225,701,680'
Set-TextCell $ws 11 5 '5.70%'
Set-TextCell $ws 11 6 '-26.44%'

Set-TextCell $ws 12 1 'Universities'
Set-TextCell $ws 12 2 '80.00%'
Set-TextCell $ws 12 3 '5'
Set-TextCell $ws 12 4 '$32,373,104'
Set-TextCell $ws 12 5 '8.81%'
Set-TextCell $ws 12 6 '-10.69%'

Set-TextCell $ws 13 1 'Total'
Set-TextCell $ws 13 2 '70.00%'
Set-TextCell $ws 13 3 '630'
Set-TextCell $ws 13 4 '$1,081,113,515'
Set-TextCell $ws 13 5 '8.92%'
Set-TextCell $ws 13 6 '-17.37%'

